$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns O, K, J (rightmost first so earlier column indices stay valid)
$ws.Columns.Item(15).Delete()   # O - Lucky
$ws.Columns.Item(11).Delete()   # K - Accuracy
$ws.Columns.Item(10).Delete()   # J - Dodge

# Rename remaining headers: I1 Defense -> Magic, O1 DefenseFloat -> MagicFloat
$ws.Range("I1").Value = "Magic"
$ws.Range("O1").Value = "MagicFloat"

# Update the selection to match the recorded state in the file
$ws.Range("R17").Select()
